$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 135
$ws.Range("H135").Value = 967
$ws.Range("I135").Value = 810.6667
$ws.Range("J135").Value = 1436
$ws.Range("K135").Value = 7296.0003
$ws.Range("L135").Value = 12924
$ws.Range("M135").Value = -4761.0003
$ws.Range("N135").Value = -17994
# Row 138
$ws.Range("H138").Value = 2764.1
$ws.Range("I138").Value = 1145.174
$ws.Range("J138").Value = 4143.185
$ws.Range("K138").Value = 3435.522
$ws.Range("L138").Value = 12429.555
$ws.Range("M138").Value = 1704.478
$ws.Range("N138").Value = -22709.555

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 15630408
$ws.Range("I32").Value = 22729336
$ws.Range("J32").Value = 12767.7
$ws.Range("K32").Value = 22729336
$ws.Range("L32").Value = 12767.7
$ws.Range("M32").Value = -22729049
$ws.Range("N32").Value = -13341.7
# Row 61
$ws.Range("H61").Value = 2073.3333
$ws.Range("I61").Value = 1953.6471
$ws.Range("J61").Value = 2582
$ws.Range("K61").Value = 1953.6471
$ws.Range("L61").Value = 2582
$ws.Range("M61").Value = -1741.6471
$ws.Range("N61").Value = -3006
# Row 74
$ws.Range("H74").Value = 3511.8647
$ws.Range("I74").Value = 4698.6665
$ws.Range("J74").Value = 1320.8462
$ws.Range("K74").Value = 4698.6665
$ws.Range("L74").Value = 1320.8462
$ws.Range("M74").Value = -3824.6665
$ws.Range("N74").Value = -3068.8462
# Row 77
$ws.Range("H77").Value = 3511.8647
$ws.Range("I77").Value = 4698.6665
$ws.Range("J77").Value = 1320.8462
$ws.Range("K77").Value = 23493.3325
$ws.Range("L77").Value = 6604.231
$ws.Range("M77").Value = -19125.3325
$ws.Range("N77").Value = -15340.231
# Row 136
$ws.Range("H136").Value = 2073.3333
$ws.Range("I136").Value = 1953.6471
$ws.Range("J136").Value = 2582
$ws.Range("K136").Value = 5860.9413
$ws.Range("L136").Value = 7746
$ws.Range("M136").Value = -3310.9413
$ws.Range("N136").Value = -12846

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 3571.4285
$ws.Range("I5").Value = 1125
$ws.Range("J5").Value = 6833.3335
$ws.Range("K5").Value = 1125
$ws.Range("L5").Value = 6833.3335
$ws.Range("M5").Value = -1012
$ws.Range("N5").Value = -7059.3335

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 1302918.1
$ws.Range("I6").Value = 8750498
$ws.Range("J6").Value = 7686.9565
$ws.Range("K6").Value = 8750498
$ws.Range("L6").Value = 7686.9565
$ws.Range("M6").Value = -8750385
$ws.Range("N6").Value = -7912.9565
# Row 16
$ws.Range("H16").Value = 3497.5557
$ws.Range("I16").Value = 1636.3636
$ws.Range("J16").Value = 6422.2856
$ws.Range("K16").Value = 1636.3636
$ws.Range("L16").Value = 6422.2856
$ws.Range("M16").Value = -1349.3636
$ws.Range("N16").Value = -6996.2856
# Row 31
$ws.Range("H31").Value = 1483.5714
$ws.Range("I31").Value = 1595.5294
$ws.Range("J31").Value = 1377.8334
$ws.Range("K31").Value = 1595.5294
$ws.Range("L31").Value = 1377.8334
$ws.Range("M31").Value = -1300.5294
$ws.Range("N31").Value = -1967.8334
# Row 34
$ws.Range("H34").Value = 1483.5714
$ws.Range("I34").Value = 1595.5294
$ws.Range("J34").Value = 1377.8334
$ws.Range("K34").Value = 1595.5294
$ws.Range("L34").Value = 1377.8334
$ws.Range("M34").Value = -1393.5294
$ws.Range("N34").Value = -1781.8334
# Row 35
$ws.Range("H35").Value = 5316.4443
$ws.Range("I35").Value = 824.6667
$ws.Range("J35").Value = 14300
$ws.Range("K35").Value = 824.6667
$ws.Range("L35").Value = 14300
$ws.Range("M35").Value = -530.6667
$ws.Range("N35").Value = -14888
# Row 86
$ws.Range("H86").Value = 3892.182
$ws.Range("I86").Value = 4312.6665
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 4312.6665
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -3189.6665
$ws.Range("N86").Value = -4246
# Row 89
$ws.Range("H89").Value = 3892.182
$ws.Range("I89").Value = 4312.6665
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 21563.3325
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -15947.3325
$ws.Range("N89").Value = -21232
# Row 94
$ws.Range("H94").Value = 3142.56
$ws.Range("I94").Value = 1398.5
$ws.Range("J94").Value = 4305.2666
$ws.Range("K94").Value = 1398.5
$ws.Range("L94").Value = 4305.2666
$ws.Range("M94").Value = -947.5
$ws.Range("N94").Value = -5207.2666
# Row 113
$ws.Range("H113").Value = 3497.5557
$ws.Range("I113").Value = 1636.3636
$ws.Range("J113").Value = 6422.2856
$ws.Range("K113").Value = 1636.3636
$ws.Range("L113").Value = 6422.2856
$ws.Range("M113").Value = 533.6364000000001
$ws.Range("N113").Value = -10762.2856
# Row 122
$ws.Range("H122").Value = 873.5714
$ws.Range("I122").Value = 1006
$ws.Range("J122").Value = 820.6
$ws.Range("K122").Value = 3018
$ws.Range("L122").Value = 2461.8
$ws.Range("M122").Value = -568
$ws.Range("N122").Value = -7361.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 457.6111
$ws.Range("I5").Value = 305.81818
$ws.Range("J5").Value = 524.4
$ws.Range("K5").Value = 917.45454
$ws.Range("L5").Value = 1573.2
$ws.Range("M5").Value = -805.45454
$ws.Range("N5").Value = -1797.2
# Row 68
$ws.Range("H68").Value = 765.7536
$ws.Range("J68").Value = 835.9375
$ws.Range("L68").Value = 2507.8125
$ws.Range("N68").Value = -4129.8125
# Row 71
$ws.Range("H71").Value = 765.7536
$ws.Range("J71").Value = 835.9375
$ws.Range("L71").Value = 7523.4375
$ws.Range("N71").Value = -15635.4375
# Row 107
$ws.Range("H107").Value = 35714816
$ws.Range("I107").Value = 255.90475
$ws.Range("K107").Value = 767.71425
$ws.Range("M107").Value = 1152.28575
# Row 131
$ws.Range("H131").Value = 2293.4866
$ws.Range("I131").Value = 741.2
$ws.Range("J131").Value = 2405.971
$ws.Range("K131").Value = 2223.6
$ws.Range("L131").Value = 7217.913
$ws.Range("M131").Value = 2816.4
$ws.Range("N131").Value = -17297.913
# Row 132
$ws.Range("H132").Value = 1123260.5
$ws.Range("I132").Value = 700
$ws.Range("J132").Value = 1443992
$ws.Range("K132").Value = 6300
$ws.Range("L132").Value = 12995928
$ws.Range("M132").Value = -3770
$ws.Range("N132").Value = -13000988
# Row 135
$ws.Range("H135").Value = 457.6111
$ws.Range("I135").Value = 305.81818
$ws.Range("J135").Value = 524.4
$ws.Range("K135").Value = 2752.36362
$ws.Range("L135").Value = 4719.599999999999
$ws.Range("M135").Value = -217.3636200000001
$ws.Range("N135").Value = -9789.599999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 41
$ws.Range("H41").Value = 13758.667
$ws.Range("I41").Value = 2499
$ws.Range("K41").Value = 2499
$ws.Range("M41").Value = -2144
# Row 122
$ws.Range("H122").Value = 3255.8
$ws.Range("I122").Value = 3453.3635
$ws.Range("J122").Value = 2712.5
$ws.Range("K122").Value = 10360.0905
$ws.Range("L122").Value = 8137.5
$ws.Range("M122").Value = -7910.0905
$ws.Range("N122").Value = -13037.5
# Row 132
$ws.Range("H132").Value = 3525.0667
$ws.Range("I132").Value = 3431.375
$ws.Range("J132").Value = 3899.8333
$ws.Range("K132").Value = 10294.125
$ws.Range("L132").Value = 11699.4999
$ws.Range("M132").Value = -7764.125
$ws.Range("N132").Value = -16759.4999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
